$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of row 6 and row 7 (columns A through AY) -
# the two species records were reordered relative to each other.

$lastCol = 51  # AY is the 51st column

for ($c = 1; $c -le $lastCol; $c++) {
    $cell6 = $ws.Cells.Item(6, $c)
    $cell7 = $ws.Cells.Item(7, $c)

    $v6 = $cell6.Value2
    $v7 = $cell7.Value2

    # Only touch cells whose value actually changes - writing an empty
    # string into an already-empty cell would delete it outright, which
    # would incorrectly remove placeholder cells that must stay in place.
    if ($v6 -ne $v7) {
        $cell6.Value2 = $v7
        $cell7.Value2 = $v6
    }
}
